$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2023 data column (T) to the right of the existing 2022 column (S),
# reusing S's per-row formatting (header style, body style, bottom-border
# style on the last data row) for the new column.
$ws.Range("S4:S14").Copy() | Out-Null
$ws.Range("T4:T14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("T4").Value = 2023

$ws.Range("T5").Value = 99.5
$ws.Range("T6").Value = 99.426175237254469
$ws.Range("T7").Value = 99.458151211935132
$ws.Range("T8").Value = 99.44178628389156
$ws.Range("T9").Value = 99.453125
$ws.Range("T10").Value = 99.487295483676391
$ws.Range("T11").Value = 99.743589743589752
$ws.Range("T12").Value = 99.190647482014398
$ws.Range("T13").Value = 99.483321247280642
$ws.Range("T14").Value = 99.771121504627331

$ws.Range("A1").Select() | Out-Null
